{"js": "// Replace each old math-problem answer text with its new value.\n// Each (old, new) pair corresponds to one table-cell run in the document;\n// all old values are unique within the document, so an exact-text search\n// unambiguously targets the single matching run.\nconst replacements = [\n  [\"0+58=58\", \"56-34=22\"],\n  [\"6+4=10\", \"9+8=17\"],\n  [\"93-43=50\", \"71-9=62\"],\n  [\"45+18=63\", \"66-40=26\"],\n  [\"90-82=8\", \"64-47=17\"],\n  [\"69-31=38\", \"29-21=8\"],\n  [\"77-72=5\", \"99-46=53\"],\n  [\"92-54=38\", \"8+0=8\"],\n  [\"49-43=6\", \"88-52=36\"],\n  [\"92-71=21\", \"18+64=82\"],\n  [\"48+43=91\", \"15+45=60\"],\n  [\"74-48=26\", \"47-20=27\"],\n  [\"53-48=5\", \"2+67=69\"],\n  [\"26+14=40\", \"48-23=25\"],\n  [\"28+68=96\", \"50+17=67\"],\n  [\"88-58=30\", \"20-1=19\"],\n  [\"68+0=68\", \"16+57=73\"],\n  [\"53-25=28\", \"18-5=13\"],\n  [\"27-9=18\", \"44+18=62\"],\n  [\"36+3=39\", \"23+50=73\"],\n  [\"29+1=30\", \"75+2=77\"],\n  [\"63+13=76\", \"26+52=78\"],\n  [\"28-23=5\", \"21-9=12\"],\n  [\"52+20=72\", \"65-29=36\"],\n  [\"14-2=12\", \"26+21=47\"],\n  [\"71+4=75\", \"92-91=1\"],\n  [\"38-23=15\", \"0+55=55\"],\n  [\"62+33=95\", \"46+5=51\"],\n  [\"80-48=32\", \"15+80=95\"],\n  [\"30+6=36\", \"25+38=63\"],\n  [\"18+58=76\", \"99-54=45\"],\n  [\"55+19=74\", \"5+24=29\"],\n  [\"28+42=70\", \"15+18=33\"],\n  [\"4+78=82\", \"26-17=9\"],\n  [\"93-6=87\", \"3+37=40\"],\n  [\"66-5=61\", \"35-0=35\"],\n  [\"9+55=64\", \"3+54=57\"],\n  [\"88-14=74\", \"45+19=64\"],\n  [\"21+45=66\", \"11+2=13\"],\n  [\"14+83=97\", \"62+28=90\"],\n  [\"35+63=98\", \"1+29=30\"],\n  [\"71-66=5\", \"17+43=60\"],\n  [\"8+40=48\", \"61-42=19\"],\n  [\"30+47=77\", \"41-1=40\"],\n  [\"30+43=73\", \"97-96=1\"],\n  [\"27+59=86\", \"80-74=6\"],\n  [\"49-5=44\", \"47+52=99\"],\n  [\"61-43=18\", \"27-15=12\"],\n  [\"7+44=51\", \"82-76=6\"],\n  [\"76-63=13\", \"91-21=70\"],\n  [\"71-39=32\", \"89-84=5\"],\n  [\"74-37=37\", \"84-38=46\"],\n  [\"16+36=52\", \"38-13=25\"],\n  [\"75-62=13\", \"75-17=58\"],\n  [\"74-12=62\", \"98-39=59\"],\n  [\"64-60=4\", \"21+17=38\"],\n  [\"71-54=17\", \"30+50=80\"],\n  [\"92-78=14\", \"66-21=45\"],\n  [\"36-25=11\", \"16+58=74\"],\n  [\"40+33=73\", \"30+29=59\"],\n  [\"65+33=98\", \"55-13=42\"],\n  [\"98-76=22\", \"18+48=66\"],\n  [\"87-70=17\", \"4+26=30\"],\n  [\"91-3=88\", \"71-21=50\"],\n  [\"84+13=97\", \"32+47=79\"],\n  [\"75-6=69\", \"84-72=12\"],\n  [\"19+8=27\", \"6+69=75\"],\n  [\"53+24=77\", \"63+25=88\"],\n  [\"82-27=55\", \"69-39=30\"],\n  [\"58-36=22\", \"31+18=49\"],\n  [\"37+52=89\", \"63-53=10\"],\n  [\"2+11=13\", \"14-10=4\"],\n  [\"92+4=96\", \"97-85=12\"],\n  [\"71-28=43\", \"39+11=50\"],\n  [\"19-10=9\", \"98-57=41\"],\n  [\"55+40=95\", \"53-31=22\"],\n  [\"68-3=65\", \"4+33=37\"],\n  [\"85-33=52\", \"56+1=57\"],\n  [\"32+35=67\", \"71-0=71\"],\n  [\"79-43=36\", \"82-44=38\"],\n  [\"96-9=87\", \"78-10=68\"],\n  [\"72-42=30\", \"59-29=30\"],\n  [\"98-5=93\", \"14+10=24\"],\n  [\"66+31=97\", \"85-23=62\"],\n  [\"46-0=46\", \"9+76=85\"],\n  [\"65-38=27\", \"11+3=14\"],\n  [\"90-63=27\", \"48+7=55\"],\n  [\"1+69=70\", \"39-10=29\"],\n  [\"57+10=67\", \"88-42=46\"],\n  [\"8+85=93\", \"64+6=70\"],\n  [\"80-57=23\", \"9+49=58\"],\n  [\"61+27=88\", \"29+52=81\"],\n  [\"43-16=27\", \"92-14=78\"],\n  [\"2+20=22\", \"38+33=71\"],\n  [\"30+14=44\", \"49-23=26\"],\n  [\"32-7=25\", \"13+42=55\"],\n  [\"15+4=19\", \"10+41=51\"],\n  [\"61-21=40\", \"34+40=74\"],\n  [\"95-1=94\", \"91-12=79\"],\n  [\"83-12=71\", \"41+7=48\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update each math-problem answer cell from its old text to the new text.\n# Every 'old' string is unique within the document, so Find/Replace against\n# the whole document body (Content) unambiguously targets a single run.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('0+58=58', '56-34=22'),\n    @('6+4=10', '9+8=17'),\n    @('93-43=50', '71-9=62'),\n    @('45+18=63', '66-40=26'),\n    @('90-82=8', '64-47=17'),\n    @('69-31=38', '29-21=8'),\n    @('77-72=5', '99-46=53'),\n    @('92-54=38', '8+0=8'),\n    @('49-43=6', '88-52=36'),\n    @('92-71=21', '18+64=82'),\n    @('48+43=91', '15+45=60'),\n    @('74-48=26', '47-20=27'),\n    @('53-48=5', '2+67=69'),\n    @('26+14=40', '48-23=25'),\n    @('28+68=96', '50+17=67'),\n    @('88-58=30', '20-1=19'),\n    @('68+0=68', '16+57=73'),\n    @('53-25=28', '18-5=13'),\n    @('27-9=18', '44+18=62'),\n    @('36+3=39', '23+50=73'),\n    @('29+1=30', '75+2=77'),\n    @('63+13=76', '26+52=78'),\n    @('28-23=5', '21-9=12'),\n    @('52+20=72', '65-29=36'),\n    @('14-2=12', '26+21=47'),\n    @('71+4=75', '92-91=1'),\n    @('38-23=15', '0+55=55'),\n    @('62+33=95', '46+5=51'),\n    @('80-48=32', '15+80=95'),\n    @('30+6=36', '25+38=63'),\n    @('18+58=76', '99-54=45'),\n    @('55+19=74', '5+24=29'),\n    @('28+42=70', '15+18=33'),\n    @('4+78=82', '26-17=9'),\n    @('93-6=87', '3+37=40'),\n    @('66-5=61', '35-0=35'),\n    @('9+55=64', '3+54=57'),\n    @('88-14=74', '45+19=64'),\n    @('21+45=66', '11+2=13'),\n    @('14+83=97', '62+28=90'),\n    @('35+63=98', '1+29=30'),\n    @('71-66=5', '17+43=60'),\n    @('8+40=48', '61-42=19'),\n    @('30+47=77', '41-1=40'),\n    @('30+43=73', '97-96=1'),\n    @('27+59=86', '80-74=6'),\n    @('49-5=44', '47+52=99'),\n    @('61-43=18', '27-15=12'),\n    @('7+44=51', '82-76=6'),\n    @('76-63=13', '91-21=70'),\n    @('71-39=32', '89-84=5'),\n    @('74-37=37', '84-38=46'),\n    @('16+36=52', '38-13=25'),\n    @('75-62=13', '75-17=58'),\n    @('74-12=62', '98-39=59'),\n    @('64-60=4', '21+17=38'),\n    @('71-54=17', '30+50=80'),\n    @('92-78=14', '66-21=45'),\n    @('36-25=11', '16+58=74'),\n    @('40+33=73', '30+29=59'),\n    @('65+33=98', '55-13=42'),\n    @('98-76=22', '18+48=66'),\n    @('87-70=17', '4+26=30'),\n    @('91-3=88', '71-21=50'),\n    @('84+13=97', '32+47=79'),\n    @('75-6=69', '84-72=12'),\n    @('19+8=27', '6+69=75'),\n    @('53+24=77', '63+25=88'),\n    @('82-27=55', '69-39=30'),\n    @('58-36=22', '31+18=49'),\n    @('37+52=89', '63-53=10'),\n    @('2+11=13', '14-10=4'),\n    @('92+4=96', '97-85=12'),\n    @('71-28=43', '39+11=50'),\n    @('19-10=9', '98-57=41'),\n    @('55+40=95', '53-31=22'),\n    @('68-3=65', '4+33=37'),\n    @('85-33=52', '56+1=57'),\n    @('32+35=67', '71-0=71'),\n    @('79-43=36', '82-44=38'),\n    @('96-9=87', '78-10=68'),\n    @('72-42=30', '59-29=30'),\n    @('98-5=93', '14+10=24'),\n    @('66+31=97', '85-23=62'),\n    @('46-0=46', '9+76=85'),\n    @('65-38=27', '11+3=14'),\n    @('90-63=27', '48+7=55'),\n    @('1+69=70', '39-10=29'),\n    @('57+10=67', '88-42=46'),\n    @('8+85=93', '64+6=70'),\n    @('80-57=23', '9+49=58'),\n    @('61+27=88', '29+52=81'),\n    @('43-16=27', '92-14=78'),\n    @('2+20=22', '38+33=71'),\n    @('30+14=44', '49-23=26'),\n    @('32-7=25', '13+42=55'),\n    @('15+4=19', '10+41=51'),\n    @('61-21=40', '34+40=74'),\n    @('95-1=94', '91-12=79'),\n    @('83-12=71', '41+7=48'),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute([ref]$oldText, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$newText, 2) | Out-Null\n}\n"}
